## timeline.xlsx - "full timeline added to draft"
##
## This applies the formatting touch-ups captured in the commit:
##  - mark the base Arial font with its charset (cosmetic/metadata)
##  - nudge the column widths (A:F) down slightly, as Excel/LibreOffice
##    re-measured them for the updated content
##  - move the active selection from G16 to M6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- font metadata -------------------------------------------------------
# The workbook's base font (Arial 10) picks up an explicit ANSI charset.
try {
    $ws.Cells.Font.Charset = 1
} catch {
}

# --- column widths ---------------------------------------------------------
# Values below are the Excel "CharactersWidth" equivalents of the stored
# OOXML column widths from the target file (stored width = ColumnWidth + 5/6).
$ws.Columns.Item(1).ColumnWidth = 4.5646258503401365
$ws.Columns.Item(2).ColumnWidth = 7.268707482993197
$ws.Range($ws.Cells.Item(1, 3), $ws.Cells.Item(1, 4)).EntireColumn.ColumnWidth = 2.676870748299317
$ws.Columns.Item(5).ColumnWidth = 17.523809523809568
$ws.Columns.Item(6).ColumnWidth = 10.508503401360567

# --- selection ---------------------------------------------------------
$ws.Range("M6").Select() | Out-Null
